$d = $word.ActiveDocument

# Locate the list paragraph that currently reads "I have made links
# clickable; those links take you to external websites outside app."
# (the original text of the last bullet in the "Extras" list).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "I have made links clickable*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'I have made links clickable' paragraph"
}

$p = $d.Paragraphs.Item($targetIndex)

# Replace the paragraph's whole text (this keeps the paragraph's pPr and
# the trailing _GoBack bookmark intact, re-anchored after the new text)
# with the new "champions league" bug-report text.
$newFirstText = "I got my champions league logo to bounce when going on to the Champions League Info page the first time, but after switching it does not work, which is something I have never found out."
$r = $d.Range($p.Range.Start, $p.Range.End)
$r.Find.Execute("I have made links clickable; those links take you to external websites outside app.", $true, $false, $false, $false, $false, $true, 1, $false, $newFirstText, 2)

# Insert a brand-new list paragraph right after it (inherits the same
# list/ListParagraph formatting) and give it the original "links
# clickable" sentence that used to live in the paragraph above.
$p2 = $d.Paragraphs.Item($targetIndex)
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item($targetIndex + 1)
$originalText = "I have made links clickable; those links take you to external websites outside app."
$p3.Range.InsertAfter($originalText)
